# Add a new survey wave "2. 3. 2021" as the last column on both sheets,
# and bump the "aktualizace" date in the two title/footer cells.

$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("data")
$wsCnt  = $wb.Worksheets.Item("pocetR")

$newDateLabel = "2. 3. 2021"

# ---------------------------------------------------------------------
# Sheet "data": percentages. New column is Z (26), header row is row 1,
# data rows are 2..76, row 77 holds the title only (column A).
# ---------------------------------------------------------------------
# Copy the header formatting (bold, centered, bordered) from the
# previous header cell (Y1) before writing the new label.
$wsData.Cells.Item(1, 25).Copy($wsData.Cells.Item(1, 26))
$wsData.Cells.Item(1, 26).Value = $newDateLabel

$dataValues = @(0.17,0.33,0.5,0.23,0.38,0.39,0.17,0.33,0.5,0.15,0.3,0.55,0.16,0.31,0.53,0.18,0.32,0.5,0.17,0.36,0.47,0.23,0.35,0.42,0.12,0.31,0.57,0.19,0.32,0.49,0.17,0.33,0.5,0.17,0.33,0.5,0.13,0.36,0.51,0.04,0.25,0.71,0.16,0.3,0.54,0.16,0.36,0.48,0.28,0.36,0.36,0.2,0.35,0.45,0.19,0.32,0.49,0.28,0.41,0.31,0.21,0.35,0.44,0.22,0.31,0.47,0.18,0.33,0.49,0.11,0.3,0.59,0.09,0.3,0.61)

for ($i = 0; $i -lt $dataValues.Length; $i++) {
    $wsData.Cells.Item($i + 2, 26).Value = $dataValues[$i]
}

$wsData.Cells.Item(77, 1).Value = "Život během pandemie, Obavy z epidemie, % respondentů celkově a ve skupinách, aktualizace 9. 3. 2021"

# ---------------------------------------------------------------------
# Sheet "pocetR": sample sizes. New column is Y (25), header row is
# row 1, data rows are 2..26, row 27 holds the title (column A) plus
# blank cells across the rest of the row.
# ---------------------------------------------------------------------
$wsCnt.Cells.Item(1, 24).Copy($wsCnt.Cells.Item(1, 25))
$wsCnt.Cells.Item(1, 25).Value = $newDateLabel

$cntValues = @(2130,512,776,842,691,727,712,1042,1088,1107,490,252,281,53,159,109,18,288,574,266,401,375,246,375,467)

for ($i = 0; $i -lt $cntValues.Length; $i++) {
    $wsCnt.Cells.Item($i + 2, 25).Value = $cntValues[$i]
}

$wsCnt.Cells.Item(27, 1).Value = "Život během pandemie, Obavy z epidemie, velikost dotázaného souboru celkově a ve skupinách, aktualizace 9. 3. 2021"
